$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.069019120957484
$ws.Range("D2").Value = 1.077823153059952
$ws.Range("E2").Value = 1.073285255758532
$ws.Range("F2").Value = 1.082915844690789
$ws.Range("I2").Value = 1.028826599955195
$ws.Range("J2").Value = 1.073955974649111
$ws.Range("K2").Value = 1.080503440761082
$ws.Range("L2").Value = 1.075977521907853
$ws.Range("M2").Value = 1.08558281856547
$ws.Range("N2").Value = 1.028199445568428

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.07143628937875
$ws.Range("D3").Value = 1.080193302035148
$ws.Range("E3").Value = 1.075484881710314
$ws.Range("F3").Value = 1.085358401357407
$ws.Range("I3").Value = 1.02908509060353
$ws.Range("J3").Value = 1.076023245003704
$ws.Range("K3").Value = 1.082687575346458
$ws.Range("L3").Value = 1.077990655524713
$ws.Range("M3").Value = 1.087840182043569
$ws.Range("N3").Value = 1.028887996310257

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.072984117457136
$ws.Range("D4").Value = 1.081709492391961
$ws.Range("E4").Value = 1.076892050579358
$ws.Range("F4").Value = 1.086917309093332
$ws.Range("I4").Value = 1.029242118693787
$ws.Range("J4").Value = 1.077344672200696
$ws.Range("K4").Value = 1.084083105105078
$ws.Range("L4").Value = 1.07927681974145
$ws.Range("M4").Value = 1.089278980483433
$ws.Range("N4").Value = 1.029328026116389

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.073631013351437
$ws.Range("D5").Value = 1.082342797843427
$ws.Range("E5").Value = 1.077479832957073
$ws.Range("F5").Value = 1.087567591488412
$ws.Range("I5").Value = 1.029305697846106
$ws.Range("J5").Value = 1.077896383097141
$ws.Range("K5").Value = 1.084665610480825
$ws.Range("L5").Value = 1.079813649395276
$ws.Range("M5").Value = 1.089878698913228
$ws.Range("N5").Value = 1.029511719154597

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.073739408729756
$ws.Range("D6").Value = 1.082448894334402
$ws.Range("E6").Value = 1.077578303904873
$ws.Range("F6").Value = 1.087676481168919
$ws.Range("I6").Value = 1.029316230643029
$ws.Range("J6").Value = 1.077988795920277
$ws.Range("K6").Value = 1.08476317296032
$ws.Range("L6").Value = 1.079903560240943
$ws.Range("M6").Value = 1.089979094468624
$ws.Range("N6").Value = 1.02954248671559

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.072992776185234
$ws.Range("D7").Value = 1.081717970662658
$ws.Range("E7").Value = 1.076899919360075
$ws.Range("F7").Value = 1.086926018056107
$ws.Range("I7").Value = 1.029242977792664
$ws.Range("J7").Value = 1.077352059096119
$ws.Range("K7").Value = 1.084090904879213
$ws.Range("L7").Value = 1.079284008016959
$ws.Range("M7").Value = 1.089287014088058
$ws.Range("N7").Value = 1.029330485691656

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.069839424684359
$ws.Range("D8").Value = 1.078627817819899
$ws.Range("E8").Value = 1.074032014784045
$ws.Range("F8").Value = 1.08374584151299
$ws.Range("I8").Value = 1.028916085755686
$ws.Range("J8").Value = 1.074658022678529
$ws.Range("K8").Value = 1.081245300549087
$ws.Range("L8").Value = 1.076661321837036
$ws.Range("M8").Value = 1.086350283557198
$ws.Range("N8").Value = 1.02843329940117

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.064154813679999
$ws.Range("D9").Value = 1.073045252643525
$ws.Range("E9").Value = 1.068851459266332
$ws.Range("F9").Value = 1.077972560265762
$ws.Range("I9").Value = 1.028260991053742
$ws.Range("J9").Value = 1.069783183462584
$ws.Range("K9").Value = 1.076091538057456
$ws.Range("L9").Value = 1.071910468255929
$ws.Range("M9").Value = 1.081004034511431
$ws.Range("N9").Value = 1.026809063352123

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.060273747633182
$ws.Range("D10").Value = 1.069225900598668
$ws.Range("E10").Value = 1.065307505447957
$ws.Range("F10").Value = 1.074003875750705
$ws.Range("I10").Value = 1.027770033257069
$ws.Range("J10").Value = 1.066442732605867
$ws.Range("K10").Value = 1.072556841162668
$ws.Range("L10").Value = 1.068651559106127
$ws.Range("M10").Value = 1.077318968639873
$ws.Range("N10").Value = 1.025695549803324

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.058570307346343
$ws.Range("D11").Value = 1.067547670010849
$ws.Range("E11").Value = 1.063750375555772
$ws.Range("F11").Value = 1.072255569971811
$ws.Range("I11").Value = 1.027544317817337
$ws.Range("J11").Value = 1.06497368053575
$ws.Range("K11").Value = 1.071001630421227
$ws.Range("L11").Value = 1.067217563141581
$ws.Range("M11").Value = 1.075693274409117
$ws.Range("N11").Value = 1.025205733094287

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.057934020794381
$ws.Range("D12").Value = 1.066920518641348
$ws.Range("E12").Value = 1.063168493227142
$ws.Range("F12").Value = 1.071601562535875
$ws.Range("I12").Value = 1.027458480334124
$ws.Range("J12").Value = 1.064424511003106
$ws.Range("K12").Value = 1.07042014265738
$ws.Range("L12").Value = 1.066681378174043
$ws.Range("M12").Value = 1.0750847855689
$ws.Range("N12").Value = 1.025022609302231

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.058070669028416
$ws.Range("D13").Value = 1.06705521772904
$ws.Range("E13").Value = 1.063293468827472
$ws.Range("F13").Value = 1.07174205999492
$ws.Range("I13").Value = 1.027476983531608
$ws.Range("J13").Value = 1.064542469697131
$ws.Range("K13").Value = 1.070545048099477
$ws.Range("L13").Value = 1.066796553285047
$ws.Range("M13").Value = 1.075215520177272
$ws.Range("N13").Value = 1.025061944113274

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.058517784838481
$ws.Range("D14").Value = 1.067495907297463
$ws.Range("E14").Value = 1.063702348892069
$ws.Range("F14").Value = 1.072201604281011
$ws.Range("I14").Value = 1.027537263343928
$ws.Range("J14").Value = 1.06492835797108
$ws.Range("K14").Value = 1.070953642891979
$ws.Range("L14").Value = 1.067173314630558
$ws.Range("M14").Value = 1.075643071776895
$ws.Range("N14").Value = 1.02519062037835

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.058792793451849
$ws.Range("D15").Value = 1.067766926121094
$ws.Range("E15").Value = 1.063953807107576
$ws.Range("F15").Value = 1.072484130105434
$ws.Range("I15").Value = 1.02757413838242
$ws.Range("J15").Value = 1.065165649844133
$ws.Range("K15").Value = 1.071204883028355
$ws.Range("L15").Value = 1.067404978251716
$ws.Range("M15").Value = 1.075905882713148
$ws.Range("N15").Value = 1.025269744144686

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.060386307275033
$ws.Range("D16").Value = 1.069336755155007
$ws.Range("E16").Value = 1.065410362691005
$ws.Range("F16").Value = 1.074119265761149
$ws.Range("I16").Value = 1.027784734614245
$ws.Range("J16").Value = 1.066539743803525
$ws.Range("K16").Value = 1.0726595266503
$ws.Range("L16").Value = 1.068746238523814
$ws.Range("M16").Value = 1.077426217580263
$ws.Range("N16").Value = 1.025727893151587

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.061379665364634
$ws.Range("D17").Value = 1.070314850342039
$ws.Range("E17").Value = 1.066317906136354
$ws.Range("F17").Value = 1.075136867337404
$ws.Range("I17").Value = 1.027913304691276
$ws.Range("J17").Value = 1.067395552421106
$ws.Range("K17").Value = 1.07356530840483
$ws.Range("L17").Value = 1.069581384694722
$ws.Range("M17").Value = 1.07837175822677
$ws.Range("N17").Value = 1.026013204448257

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.061956871269242
$ws.Range("D18").Value = 1.07088300779776
$ws.Range("E18").Value = 1.066845089985466
$ws.Range("F18").Value = 1.07572754687332
$ws.Range("I18").Value = 1.027987032570145
$ws.Range("J18").Value = 1.067892556135182
$ws.Range("K18").Value = 1.074091263428414
$ws.Range("L18").Value = 1.070066311970602
$ws.Range("M18").Value = 1.078920386481344
$ws.Range("N18").Value = 1.026178885157277

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.062153313166185
$ws.Range("D19").Value = 1.071076339711131
$ws.Range("E19").Value = 1.067024480897345
$ws.Range("F19").Value = 1.075928470211696
$ws.Range("I19").Value = 1.028011958044035
$ws.Range("J19").Value = 1.06806165562403
$ws.Range("K19").Value = 1.074270201335191
$ws.Range("L19").Value = 1.070231289539305
$ws.Range("M19").Value = 1.079106967991318
$ws.Range("N19").Value = 1.026235254060558

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.061273316100088
$ws.Range("D20").Value = 1.070210153739542
$ws.Range("E20").Value = 1.066220760591858
$ws.Range("F20").Value = 1.075027986093475
$ws.Range("I20").Value = 1.027899641342426
$ws.Range("J20").Value = 1.067303957936483
$ws.Range("K20").Value = 1.073468372721619
$ws.Range("L20").Value = 1.069492009609682
$ws.Range("M20").Value = 1.078270610406651
$ws.Range("N20").Value = 1.025982669668198

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.058386219367464
$ws.Range("D21").Value = 1.067366240524267
$ws.Range("E21").Value = 1.063582041234769
$ws.Range("F21").Value = 1.072066408207493
$ws.Range("I21").Value = 1.027519567765051
$ws.Range("J21").Value = 1.064814820888958
$ws.Range("K21").Value = 1.070833428063403
$ws.Range("L21").Value = 1.067062466209912
$ws.Range("M21").Value = 1.075517297397
$ws.Range("N21").Value = 1.025152761387739

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.056550372352482
$ws.Range("D22").Value = 1.065556219911509
$ws.Range("E22").Value = 1.061902698356805
$ws.Range("F22").Value = 1.070177614871781
$ws.Range("I22").Value = 1.027269034763175
$ws.Range("J22").Value = 1.063229509145785
$ws.Range("K22").Value = 1.069154614566493
$ws.Range("L22").Value = 1.065514410761908
$ws.Range("M22").Value = 1.073759307514697
$ws.Range("N22").Value = 1.024624096838554

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.057525581796965
$ws.Range("D23").Value = 1.066517864426389
$ws.Range("E23").Value = 1.062794907381355
$ws.Range("F23").Value = 1.071181477348422
$ws.Range("I23").Value = 1.027402951918858
$ws.Range("J23").Value = 1.064071871144219
$ws.Range("K23").Value = 1.070046719281726
$ws.Range("L23").Value = 1.066337042301855
$ws.Range("M23").Value = 1.074693840312597
$ws.Range("N23").Value = 1.024905014535694

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.061321377556024
$ws.Range("D24").Value = 1.070257468879234
$ws.Range("E24").Value = 1.06626466318387
$ws.Range("F24").Value = 1.075077193702559
$ws.Range("I24").Value = 1.02790581912876
$ws.Range("J24").Value = 1.067345352256247
$ws.Range("K24").Value = 1.073512181107961
$ws.Range("L24").Value = 1.069532401153706
$ws.Range("M24").Value = 1.078316323674491
$ws.Range("N24").Value = 1.025996469294918

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.065640118342238
$ws.Range("D25").Value = 1.074505284942765
$ws.Range("E25").Value = 1.070206291007077
$ws.Range("F25").Value = 1.079485760284381
$ws.Range("I25").Value = 1.028439823532324
$ws.Range("J25").Value = 1.071059044842525
$ws.Range("K25").Value = 1.077440944476843
$ws.Range("L25").Value = 1.073154477078272
$ws.Range("M25").Value = 1.082407043118302
$ws.Range("N25").Value = 1.027234255435885
